# Add a "year" column to the course upload sheet (between "end_date" and
# "grade_add_start_date"), matching the commit "added year field in the xls file".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I; everything from the old column I onward shifts right.
$ws.Columns("I:I").Insert()

# Give the new column a custom number format (0;[Red]0) like the rest of the
# workbook's formatted columns, then fill in the header and the four year values.
$ws.Range("I1:I5").NumberFormat = "0;[Red]0"

$ws.Range("I1").Value = "year"
$ws.Range("I2").Value = 2018
$ws.Range("I3").Value = 2019
$ws.Range("I4").Value = 2020
$ws.Range("I5").Value = 2021

# Keep the page orientation explicit (portrait), as seen after the edit.
$ws.PageSetup.Orientation = 1

# Reflect the author's final selection/view on the new layout.
$ws.Range("I5").Select()
